$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header shared strings (Volume number, week-of dates) ---
$ws.Range("A8").Characters(21,2).Text = "38"
$ws.Range("C9").Characters(27,9).Text = "9/18/2023"
$ws.Range("C9").Characters(47,9).Text = "9/24/2023"

# --- Simple numeric value updates ---
$ws.Range("M15").Value = -33.333333333333
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -42.857142857142
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 128
$ws.Range("J16").Value = 132
$ws.Range("K16").Value = -3.030303030303
$ws.Range("L16").Value = 5.785123966942
$ws.Range("M16").Value = -22.424242424242
$ws.Range("N16").Value = -79.905808477237
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 44.444444444444
$ws.Range("I17").Value = 186
$ws.Range("J17").Value = 166
$ws.Range("K17").Value = 12.048192771084
$ws.Range("L17").Value = -1.587301587301
$ws.Range("M17").Value = 30.985915492957
$ws.Range("N17").Value = -65.682656826568
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -60
$ws.Range("J18").Value = 131
$ws.Range("K18").Value = -35.114503816793
$ws.Range("L18").Value = -27.350427350427
$ws.Range("M18").Value = 2.409638554216
$ws.Range("N18").Value = -86.614173228346
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -32.432432432432
$ws.Range("I19").Value = 275
$ws.Range("J19").Value = 298
$ws.Range("K19").Value = -7.718120805369
$ws.Range("L19").Value = 0.7326007326
$ws.Range("M19").Value = 23.318385650224
$ws.Range("N19").Value = -20.058139534883
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 49
$ws.Range("J20").Value = 38
$ws.Range("K20").Value = 28.947368421052
$ws.Range("L20").Value = 32.432432432432
$ws.Range("M20").Value = 206.25
$ws.Range("N20").Value = -53.333333333333
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -42.307692307692
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = -12.765957446808
$ws.Range("I21").Value = 732
$ws.Range("J21").Value = 775
$ws.Range("K21").Value = -5.548387096774
$ws.Range("L21").Value = -1.744966442953
$ws.Range("M21").Value = 13.841368584758
$ws.Range("N21").Value = -68.407423392317
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = 66.666666666666
$ws.Range("I22").Value = 17
$ws.Range("K22").Value = -39.285714285714
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 112.5
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = -30
$ws.Range("I23").Value = 65
$ws.Range("J23").Value = 72
$ws.Range("K23").Value = -9.722222222222
$ws.Range("L23").Value = 3.174603174603
$ws.Range("M23").Value = 103.125
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -56
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -28.431372549019
$ws.Range("I24").Value = 940
$ws.Range("J24").Value = 862
$ws.Range("K24").Value = 9.048723897911
$ws.Range("L24").Value = 1.511879049676
$ws.Range("M24").Value = 31.837307152875
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 18.518518518518
$ws.Range("I25").Value = 304
$ws.Range("J25").Value = 313
$ws.Range("K25").Value = -2.875399361022
$ws.Range("L25").Value = 2.013422818791
$ws.Range("M25").Value = -18.933333333333
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 39
$ws.Range("K27").Value = -28.205128205128
$ws.Range("L27").Value = -22.222222222222
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -35.714285714285
$ws.Range("L28").Value = -55
$ws.Range("N28").Value = -88.461538461538
$ws.Range("J29").Value = 14
$ws.Range("K29").Value = -42.857142857142
$ws.Range("L29").Value = -57.894736842105
$ws.Range("N29").Value = -88.888888888888

# --- Cells changing between numeric and text ("0" / "***.*") representations ---
# Strategy: force text via NumberFormat "@" then assign string (for num->text),
# or assign numeric Value directly (for text->num); then PasteSpecial(formats only)
# from a same-column reference cell that already carries the exact target style.
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C20").Value = 1
$ws.Range("D27").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 1
$ws.Range("D27").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = 0
$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("D27").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("G28").Value = 1
$ws.Range("D27").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("D27").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1
$ws.Range("D27").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H29").PasteSpecial(-4122)
